# "Input Data Base Case" - bump generator P_max current (MW) / MaxInv (MW)
# from 200 to 400 for every generator row (2-17) on the Generators_New
# sheet, then leave the Python_Gen_N_Data sheet active/selected (it picks
# up the new 400 values automatically through its cross-sheet formulas).

$wb = $excel.ActiveWorkbook

$wsGen = $wb.Worksheets.Item("Generators_New")
$wsPy  = $wb.Worksheets.Item("Python_Gen_N_Data")

# Column I ("MaxInv (MW)") literal values: 200 -> 400
$wsGen.Range("I2:I17").Value = 400

# Column D ("P_max (MW)") used to hold "=I<row>" formulas that mirrored
# column I; replace them with the plain literal value 400 (formula removed).
$wsGen.Range("D2:D17").Value = 400

# Selection / active-sheet bookkeeping to match the saved view state.
$wsGen.Range("D2:D17").Select() | Out-Null

$wsPy.Activate() | Out-Null
$wsPy.Range("H19").Select() | Out-Null
